{"js": "const pairs = [\n  [\"2024-12-24 Tuesday\", \"2024-12-25 Wednesday\"],\n  [\"50\u00f73=16, 2\", \"27\u00f78=3, 3\"],\n  [\"48\u00f75=9, 3\", \"78\u00f75=15, 3\"],\n  [\"87\u00f78=10, 7\", \"55\u00f75=11, 0\"],\n  [\"31\u00f72=15, 1\", \"84\u00f72=42, 0\"],\n  [\"69\u00f74=17, 1\", \"47\u00f75=9, 2\"],\n  [\"62\u00f79=6, 8\", \"51\u00f78=6, 3\"],\n  [\"45\u00f78=5, 5\", \"67\u00f79=7, 4\"],\n  [\"33\u00f75=6, 3\", \"28\u00f77=4, 0\"],\n  [\"25\u00f75=5, 0\", \"24\u00f72=12, 0\"],\n  [\"22\u00f77=3, 1\", \"62\u00f72=31, 0\"],\n  [\"37\u00f75=7, 2\", \"22\u00f76=3, 4\"],\n  [\"51\u00f72=25, 1\", \"50\u00f77=7, 1\"],\n  [\"49\u00f75=9, 4\", \"77\u00f77=11, 0\"],\n  [\"84\u00f78=10, 4\", \"41\u00f76=6, 5\"],\n  [\"84\u00f79=9, 3\", \"13\u00f79=1, 4\"],\n  [\"38\u00f75=7, 3\", \"43\u00f76=7, 1\"],\n  [\"66\u00f72=33, 0\", \"32\u00f74=8, 0\"],\n  [\"30\u00f76=5, 0\", \"36\u00f79=4, 0\"],\n  [\"11\u00f72=5, 1\", \"24\u00f74=6, 0\"],\n  [\"80\u00f74=20, 0\", \"83\u00f78=10, 3\"],\n  [\"20\u00f79=2, 2\", \"95\u00f77=13, 4\"],\n  [\"21\u00f76=3, 3\", \"19\u00f78=2, 3\"],\n  [\"34\u00f74=8, 2\", \"34\u00f79=3, 7\"],\n  [\"29\u00f77=4, 1\", \"46\u00f74=11, 2\"],\n  [\"79\u00f78=9, 7\", \"68\u00f74=17, 0\"],\n];\n\nconst body = context.document.body;\n\nfor (const [before, after] of pairs) {\n  const results = body.search(before, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + before);\n  }\n  for (const item of results.items) {\n    item.insertText(after, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$pairs = @(\n    @(\"2024-12-24 Tuesday\", \"2024-12-25 Wednesday\"),\n    @(\"50\u00f73=16, 2\", \"27\u00f78=3, 3\"),\n    @(\"48\u00f75=9, 3\", \"78\u00f75=15, 3\"),\n    @(\"87\u00f78=10, 7\", \"55\u00f75=11, 0\"),\n    @(\"31\u00f72=15, 1\", \"84\u00f72=42, 0\"),\n    @(\"69\u00f74=17, 1\", \"47\u00f75=9, 2\"),\n    @(\"62\u00f79=6, 8\", \"51\u00f78=6, 3\"),\n    @(\"45\u00f78=5, 5\", \"67\u00f79=7, 4\"),\n    @(\"33\u00f75=6, 3\", \"28\u00f77=4, 0\"),\n    @(\"25\u00f75=5, 0\", \"24\u00f72=12, 0\"),\n    @(\"22\u00f77=3, 1\", \"62\u00f72=31, 0\"),\n    @(\"37\u00f75=7, 2\", \"22\u00f76=3, 4\"),\n    @(\"51\u00f72=25, 1\", \"50\u00f77=7, 1\"),\n    @(\"49\u00f75=9, 4\", \"77\u00f77=11, 0\"),\n    @(\"84\u00f78=10, 4\", \"41\u00f76=6, 5\"),\n    @(\"84\u00f79=9, 3\", \"13\u00f79=1, 4\"),\n    @(\"38\u00f75=7, 3\", \"43\u00f76=7, 1\"),\n    @(\"66\u00f72=33, 0\", \"32\u00f74=8, 0\"),\n    @(\"30\u00f76=5, 0\", \"36\u00f79=4, 0\"),\n    @(\"11\u00f72=5, 1\", \"24\u00f74=6, 0\"),\n    @(\"80\u00f74=20, 0\", \"83\u00f78=10, 3\"),\n    @(\"20\u00f79=2, 2\", \"95\u00f77=13, 4\"),\n    @(\"21\u00f76=3, 3\", \"19\u00f78=2, 3\"),\n    @(\"34\u00f74=8, 2\", \"34\u00f79=3, 7\"),\n    @(\"29\u00f77=4, 1\", \"46\u00f74=11, 2\"),\n    @(\"79\u00f78=9, 7\", \"68\u00f74=17, 0\"),\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $pairs) {\n    $beforeText = $pair[0]\n    $afterText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $beforeText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $afterText\n    $found = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n    if (-not $found) {\n        Write-Output \"WARNING: not found: $beforeText\"\n    }\n}\n"}
